$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tasks")

# Row 2: meet_kat
$ws.Range("A2").Value = "meet_kat"
$ws.Range("B2").Value = "Bienvenido a Celestia"
$ws.Range("C2").Value = "Has llegado a Celestia, tu nuevo hogar y lugar de trabajo. Te han hablado de Kat, una chica del lugar que te ayudara a comenzar tu nueva vida."

# Row 3: meet_chencho
$ws.Range("A3").Value = "meet_chencho"
$ws.Range("B3").Value = "Tu viejo amigo"
$ws.Range("C3").Value = "Kat ha mencionado a Chencho, un amigo de tu antigua vida que ahora vive tambien en Celestia. Toma el autobus y reunete con el en la ciudad."

# Row 4: meet_roy
$ws.Range("A4").Value = "meet_roy"
$ws.Range("B4").Value = "Routed Inc"
$ws.Range("C4").Value = "Te han hablado e Roy, uno de los empleados de Routed Inc, con quien ahora trabajaras, ve que tiene para decirte."

# Row 5: create_straight_cable
$ws.Range("A5").Value = "create_straight_cable"
$ws.Range("B5").Value = "Crea un cable directo"
$ws.Range("C5").Value = "Crea un cable directo"

# Row 6: create_crossover_cable
$ws.Range("A6").Value = "create_crossover_cable"
$ws.Range("B6").Value = "Crea dos cables cruzados"
$ws.Range("C6").Value = "Crea dos cables cruzados"

# Row 7: subnetting_reception
$ws.Range("A7").Value = "subnetting_reception"
$ws.Range("B7").Value = "El primer trabajo"
$ws.Range("C7").Value = "Has hablado con Roy y te ha dado tu primer trabajo en Routed Inc. Dirigete a la recepción de los edificios de la empresa, cerca de tu casa, y habla con Ale la encargada de la recepción."
$ws.Range("E7").Value = "reception"
$ws.Range("F7").Value = "subnetting_reception"

# Row 8: mysteries_of_celestia
$ws.Range("A8").Value = "mysteries_of_celestia"
$ws.Range("B8").Value = "Los misterios de Celestia"
$ws.Range("C8").Value = "En multiples ocasiones has oido hablar de Celestia, y algunas cosas parecen envolver al pueblo en misterio. Habla con los habitantes del pueblo para conocer mas acerca del lugar y lo que puede haber detrás de el."
